# Updated symbol list on Sat Dec 24 14:34:21 UTC 2022 with GitHub Actions
#
# This script updates price values (column D) for several coins, and
# reshuffles three rows (41-43) where BKEXToken / CEJI / KickToken swap
# order and get refreshed prices/links/labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "244.52"). Plain assignment
# via .Value would make Excel auto-convert these into real numbers, so
# first force the whole price column to Text format, make the edits,
# then restore the default ("Normal") style so no stray formatting is
# left behind on the cells.
$priceRange = $ws.Range("D2:D48")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value  = "244.52"
$ws.Range("D4").Value  = "5.385"
$ws.Range("D5").Value  = "0.05988"
$ws.Range("D7").Value  = "0.8143"
$ws.Range("D8").Value  = "0.9581"
$ws.Range("D9").Value  = "0.1426"
$ws.Range("D10").Value = "0.07426"
$ws.Range("D11").Value = "0.03276"
$ws.Range("D15").Value = "0.001595"
$ws.Range("D16").Value = "0.04810"
$ws.Range("D17").Value = "0.0005911"
$ws.Range("D18").Value = "0.005427"
$ws.Range("D19").Value = "0.004146"
$ws.Range("D20").Value = "0.0009904"
$ws.Range("D22").Value = "3.679"
$ws.Range("D23").Value = "6.432"
$ws.Range("D40").Value = "0.04021"

# Rows 41-43 swap coin identities (BKEXToken / CEJI / KickToken rotate)
# and receive refreshed prices, links and labels.
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006432"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1072"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002901"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "0.005729"
$ws.Range("D45").Value = "0.00005214"
$ws.Range("D47").Value = "0.8602"
$ws.Range("D48").Value = "0.006829"

# Restore plain/default formatting to the price column now that the
# text values are safely stored.
$priceRange.Style = "Normal"
